# Leave Balance Code - Custom leave cycle
# Add a new test case row (row 38) for the "LeaveBalance" custom accrual test case.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 38

$ws.Range("A$newRow").Value = "38"
$ws.Range("B$newRow").Value = "LeaveBalance"
$ws.Range("C$newRow").Value = "LeaveBalance"
$ws.Range("D$newRow").Value = "com.darwinbox.leaves.Accural.Custom.LeaveBalance"
$ws.Range("E$newRow").Value = "Accural//LeaveBalance.xlsx"
$ws.Range("F$newRow").Value = "LeaveBalance"
$ws.Range("G$newRow").Value = "All"

# Columns A and G are formatted as Text (same as the rest of the table)
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("G$newRow").NumberFormat = "@"

# Move the view / selection to the newly added row, mirroring the saved workbook state
$aw = $excel.ActiveWindow
$aw.ScrollRow = 20
$aw.ScrollColumn = 4
$ws.Range("G$newRow").Select()

$wb.Save()
